# Add a new logged working-hours entry (2014-04-15, 08:30-10:00) as the new
# row 126, pushing the old "empty entry row" and the summary rows
# (sum [min], sum [h], sum [working weeks]) down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- shift the trailing block (old rows 126-129) down by one row ------------
$ws.Rows.Item(126).Insert()

# --- fill in the new data row 126 -------------------------------------------
$ws.Range("A126").Value = 2014
$ws.Range("B126").Value = 4
$ws.Range("C126").Value = 15
$ws.Range("D126").Value = 0.35416666666666669
$ws.Range("E126").Value = 0.41666666666666669
$ws.Range("D126").NumberFormat = "hh:mm;@"
$ws.Range("E126").NumberFormat = "hh:mm;@"
$ws.Range("F126").Formula = "=(E126-D126)*24*60"
$ws.Range("G126").Formula = "=F126/60"
$ws.Range("F126").NumberFormat = "0"
$ws.Range("G126").NumberFormat = "0.00"

# --- the now-empty placeholder row (was row 126, now row 127) --------------
$ws.Range("D127").NumberFormat = "hh:mm;@"
$ws.Range("E127").NumberFormat = "hh:mm;@"
$ws.Range("F127").NumberFormat = "0"

# --- fix up the summary formulas so they reference the new layout ----------
$ws.Range("F128").Formula = "=SUM(F2:F127)"
$ws.Range("F129").Formula = "=F128/60"
$ws.Range("F130").Formula = "=F129/38.5"

# --- move the "active cell" selection to match the saved workbook ----------
$ws.Range("F126").Select()
